# Add a "UK" market test-data sheet, cloned from the existing "Poland" sheet,
# then adjust its content to match the UK data set.

$wb = $excel.ActiveWorkbook

$poland = $wb.Worksheets.Item("Poland")

# Duplicate "Poland" and place the copy immediately after it.
$poland.Copy($null, $poland)
$uk = $wb.Worksheets.Item("Poland (2)")
$uk.Name = "UK"

# Insert a new row at position 9 ("GMPIM"), pushing the existing rows
# 9-11 (PR1D2 / Wg / Miscellaneous) down to 10-12, and bring over the
# border formatting from the row right below so the new row matches its
# neighbours.
$uk.Range("A9:D9").Insert(-4121)
$uk.Range("A10").Copy()
$uk.Range("A9").PasteSpecial(-4122)
$uk.Range("A9").Value = "GMPIM"

# Fill in the UK-specific values (order matters so new shared strings end
# up appended in the same sequence as the source workbook).
$uk.Range("B4").Value = "NGC-2741/T3355/T3357/T3349"
$uk.Range("B2").Value = "UK Market"

# Match the saved selection/active-cell for the new sheet.
$uk.Range("B4").Select()
